$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet "Exam Sample A": append a new "Question 8" block (rows 53-59)
# =========================================================================
$wsA = $wb.Worksheets.Item("Exam Sample A")

# Row 53: blank separator row (mirrors the existing empty-row pattern
# used before every new question block on this sheet, e.g. row 43)
$wsA.Cells.Item(53, 1).Style = "Normal"
$wsA.Cells.Item(53, 2).Style = "Normal"
$wsA.Cells.Item(53, 3).Style = "Normal"

# Row 54: repeated "Questions" / "Answer" column header
$wsA.Cells.Item(54, 1).Style = "Normal"
$wsA.Cells.Item(54, 2).Value = "Questions"
$wsA.Cells.Item(54, 3).Value = "Answer"

$q8Label = "Question 8:"
$q8Full = "Question #8`nMatch the following test work products (1-4) with the right description (A-D).`n1. Test suite.`n2. Test case.`n3. Test script.`n4. Test charter.`nA. A set of test scripts to be executed in a specific test run.`nB. A set of instructions for the execution of a test.`nC. Contains expected results.`nD. Documentation of test activities in session-based exploratory testing."

$q8Answers = @(
    "Match the following test work products (1-4) with the right description (A-D).",
    "a) 1A, 2C, 3B, 4D.",
    "b) 1D, 2B, 3A, 4C.",
    "c) 1A, 2C, 3D, 4B.",
    "d) 1D, 2C, 3B, 4A."
)

$row = 55
foreach ($answer in $q8Answers) {
    $wsA.Cells.Item($row, 1).Value = $q8Label
    $wsA.Cells.Item($row, 2).Value = $q8Full
    $wsA.Cells.Item($row, 3).Value = $answer
    $wsA.Rows.Item($row).AutoFit()
    $row = $row + 1
}

# =========================================================================
# Sheet "Exam Sample B": extend the "Question 38" block with two more
# answer choices (rows 16-22)
# =========================================================================
$wsB = $wb.Worksheets.Item("Exam Sample B")

# Row 16: blank separator row
$wsB.Cells.Item(16, 1).Style = "Normal"
$wsB.Cells.Item(16, 2).Style = "Normal"
$wsB.Cells.Item(16, 3).Style = "Normal"

# Row 17: repeated "Questions" / "Answer" column header
$wsB.Cells.Item(17, 1).Style = "Normal"
$wsB.Cells.Item(17, 2).Value = "Questions"
$wsB.Cells.Item(17, 3).Value = "Answer"

$q38Label = "Question 38:"
$q38Full = "Question #38`nYou are performing system testing of a train reservation system. Based on the test cases performed, you have noticed that the system occasionally reports that no trains are available, although this should actually be the case. You have provided the developers with a summary of the defect and the version of the tested system. They recognize the urgency of the defect and are now waiting for you to provide further details.`nIn addition to the information already provided, the following additional information is given:`n1. Degree of impact (severity) of the defect.`n2. Identification of the test item.`n3. Details of the test environment.`n4. Urgency/priority to fix.`n5. Actual results.`n6. Reference to test case specification.`nWhich of this information is most useful to include in the defect report?"

$q38Answers = @(
    "1. Degree of impact (severity) of the defect.",
    "a) 1, 2, 6",
    "b) 1, 4, 5, 6",
    "c) 2, 3, 4, 5",
    "d) 3, 5, 6"
)

$row = 18
foreach ($answer in $q38Answers) {
    $wsB.Cells.Item($row, 1).Value = $q38Label
    $wsB.Cells.Item($row, 2).Value = $q38Full
    $wsB.Cells.Item($row, 3).Value = $answer
    $wsB.Rows.Item($row).AutoFit()
    $row = $row + 1
}
